$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 170, shifting existing rows 170-277 down to 171-278.
$ws.Rows.Item(170).Insert()

# Populate the new row 170 with a fresh record (same structure as the row that
# used to occupy position 170, now at 171), with updated date/volume/price figures.
$ws.Cells.Item(170, 1).Value()  = 3
$ws.Cells.Item(170, 2).Value()  = "Femacal de La Calera"
$ws.Cells.Item(170, 3).Value()  = "Coquimbo"
$ws.Cells.Item(170, 4).Value()  = 44596
$ws.Cells.Item(170, 5).Value()  = 5
$ws.Cells.Item(170, 6).Value()  = 100112012
$ws.Cells.Item(170, 7).Value()  = "Espinaca"
$ws.Cells.Item(170, 8).Value()  = "Sin especificar"
$ws.Cells.Item(170, 9).Value()  = "Primera"
$ws.Cells.Item(170, 10).Value() = 130
$ws.Cells.Item(170, 11).Value() = 3500
$ws.Cells.Item(170, 12).Value() = 4000
$ws.Cells.Item(170, 13).Value() = 3769
$ws.Cells.Item(170, 14).Value() = "$/docena de atados (3 kilos)"
$ws.Cells.Item(170, 15).Value() = "Provincia de Quillota"
$ws.Cells.Item(170, 16).Value() = 1256
$ws.Cells.Item(170, 17).Value() = 3
$ws.Cells.Item(170, 18).Value() = "Hortaliza"
